$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 duplicates the values already present in row 2 (columns A-O).
# Columns I-O look numeric ("5","6","2",...) but must stay text cells,
# exactly like the matching cells in row 2, so a leading apostrophe is
# used to force Excel to store them as text instead of numbers.

$ws.Cells.Item(3, 1).Value  = "slick-list\ draggable""] [class=""slide\ slick-slide"
$ws.Cells.Item(3, 2).Value  = "slick-list\ draggable""] [class=""slide\ slick-slide"
$ws.Cells.Item(3, 3).Value  = "slick-list\ draggable""] [class=""slide\ slick-slide"
$ws.Cells.Item(3, 4).Value  = "slick-list\ draggable""] [class=""slide\ slick-slide"
$ws.Cells.Item(3, 5).Value  = "slick-list\ draggable""] [class=""slide\ slick-slide"
$ws.Cells.Item(3, 6).Value  = " "
$ws.Cells.Item(3, 7).Value  = " "
$ws.Cells.Item(3, 8).Value  = " "
$ws.Cells.Item(3, 9).Value  = "'5"
$ws.Cells.Item(3, 10).Value = "'6"
$ws.Cells.Item(3, 11).Value = "'2"
$ws.Cells.Item(3, 12).Value = "'3"
$ws.Cells.Item(3, 13).Value = "'4"
$ws.Cells.Item(3, 14).Value = "'5"
$ws.Cells.Item(3, 15).Value = "'6"
